$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (LimType): "UP" -> "FX"
$ws.Range("C7").Value = "FX"
$ws.Range("C8").Value = "FX"
$ws.Range("C9").Value = "FX"
$ws.Range("C10").Value = "FX"

# Column E (Wartość atrybutu / attribute value)
$ws.Range("E7").Value = 0
$ws.Range("E8").Value = 3.33
$ws.Range("E9").Value = 3.75
$ws.Range("E10").Value = 0

# Column F (Filtr procesow: nazwa procesu): "EL_ATOMOWA" -> "EL_NEW_ATOMOWA"
$ws.Range("F7").Value = "EL_NEW_ATOMOWA"
$ws.Range("F8").Value = "EL_NEW_ATOMOWA"
$ws.Range("F9").Value = "EL_NEW_ATOMOWA"
$ws.Range("F10").Value = "EL_NEW_ATOMOWA"
